$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (existing row, values changed; sending/target cluster pair ECs->ECs)
$row2 = @("ECs","Ly9","Ly9","ECs",1,0.3333333333333333,5.025175333333333,15.075526,0.1274621509804351,0.1274621509804351,1,0.3333333333333333,5.025175333333333,15.075526,0.1274621509804351,0.1274621509804351,25.25238713074178,227.271484176676,0.01624659993255924,0.01624659993255923)

# Row 3 (new) - ECs -> M2
$row3 = @("ECs","Ly9","Ly9","M2",1,0.3333333333333333,5.025175333333333,15.075526,0.1274621509804351,0.1274621509804351,3,1,34.399668,103.199004,0.872537849019565,0.8725378490195649,172.864363108456,1555.779267976104,0.1112155510478759,0.1112155510478759)

# Row 4 (new) - M2 -> ECs
$row4 = @("M2","Ly9","Ly9","ECs",3,1,34.399668,103.199004,0.872537849019565,0.8725378490195649,1,0.3333333333333333,5.025175333333333,15.075526,0.1274621509804351,0.1274621509804351,172.864363108456,1555.779267976104,0.1112155510478759,0.1112155510478759)

# Row 5 (new) - M2 -> M2
$row5 = @("M2","Ly9","Ly9","M2",3,1,34.399668,103.199004,0.872537849019565,0.8725378490195649,3,1,34.399668,103.199004,0.872537849019565,0.8725378490195649,1183.337158510224,10650.03442659202,0.7613222979716892,0.761322297971689)

$rows = @($row2, $row3, $row4, $row5)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
